# The document originally opens with three paragraphs that get removed
# entirely in this edit:
#   1. "OVERNIGHT EVENTS: No acute events were noted over night. "
#   2. (an empty paragraph)
#   3. "SUBJECTIVE: "
# leaving the "I reviewed the following medical systems..." paragraph as
# the new first paragraph of the document body.

$d = $word.ActiveDocument

# Build a range spanning from the start of paragraph 1 through the end of
# paragraph 3 (inclusive of paragraph 3's trailing mark) and delete it in
# one shot, which removes all three paragraphs (and their marks) while
# leaving paragraph 4 ("I reviewed the following medical systems...")
# untouched as the new first paragraph.
$startRange = $d.Paragraphs.Item(1).Range.Start
$endRange = $d.Paragraphs.Item(3).Range.End
$r = $d.Range($startRange, $endRange)
$r.Delete()
